$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by duplicating "2022-Q2" (this
#    keeps identical styles/number formats) and placing it right
#    before "2022-Q2", then overwrite its data.
# ------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# extend the styled/bordered "index" column (A) down to rows 4-7
$q3.Range("A2").Copy($q3.Range("A4:A7"))

# text columns (B..G) must stay as TEXT, not be auto-converted to numbers
$q3.Range("B2:G7").NumberFormat = "@"

$idx = New-Object 'object[,]' 6,1
$idx[0,0] = 0
$idx[1,0] = 1
$idx[2,0] = 2
$idx[3,0] = 3
$idx[4,0] = 4
$idx[5,0] = 5
$q3.Range("A2:A7").Value = $idx

$data = New-Object 'object[,]' 6,6
$rows = @(
  @("016935", "景顺长城中证500指数增强C", "15.57", "93.89", "1.77", "0.2756"),
  @("008851", "景顺长城量化对冲策略三个月定期开放灵活配置混合", "2.96", "64.77", "1.18", "0.0349"),
  @("519615", "银河君尚灵活配置混合I", "3.59", "35.36", "0.73", "0.0262"),
  @("519613", "银河君尚灵活配置混合A", "2.10", "35.36", "0.73", "0.0153"),
  @("519614", "银河君尚灵活配置混合C", "0.17", "35.36", "0.73", "0.0012"),
  @("006682", "景顺长城中证500指数增强A", "0.00", "93.89", "1.77", "")
)
for ($i = 0; $i -lt $rows.Count; $i++) {
  for ($j = 0; $j -lt 6; $j++) {
    $data[$i, $j] = $rows[$i][$j]
  }
}
$q3.Range("B2:G7").Value = $data

# G7 is numeric 0 (not text) in the source data
$q3.Range("G7").NumberFormat = "General"
$q3.Range("G7").Value = 0

$rank = New-Object 'object[,]' 6,1
$rank[0,0] = 8
$rank[1,0] = 10
$rank[2,0] = 8
$rank[3,0] = 8
$rank[4,0] = 8
$rank[5,0] = 8
$q3.Range("H2:H7").Value = $rank

# ------------------------------------------------------------------
# 2. Update the "总计" summary sheet: existing row 2 becomes the new
#    2022-Q3 totals, the old row 2 (2022-Q2) shifts to row 3, and the
#    old row 3 (2021-Q1) shifts to row 4.
# ------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

# push old row 3 down to row 4 (copy whole row incl. styles)
$tot.Range("A3:D3").Copy($tot.Range("A4:D4"))
$tot.Range("A4").Value = 2

# old row 2 -> row 3 ("2022-Q2" data, unchanged values)
$tot.Range("B3").Value = "2022-Q2"
$tot.Range("C3").Value = 2
$tot.Range("D3").Value = 0.06

# row 2 becomes the brand-new "2022-Q3" totals
$tot.Range("B2").Value = "2022-Q3"
$tot.Range("C2").Value = 6
$tot.Range("D2").Value = 0.35
